# ----------------------------------------------------------------------------
# Updates the "EoCEDwEC" workbook so it reflects the newer EIA source
# ("Price Elasticities for Energy Use in Buildings of the United States",
# 2014 Appendix) instead of the old 2003 NEMS table, switches the elasticity
# baseline from the 1-year to the 3-year short-run figure, refreshes the
# EIA Table 1 data, and adds kerosene / heavy-fuel-oil / LPG / hydrogen rows
# to the EoCEDwEC fuel-elasticity table.
# ----------------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$wsAbout = $wb.Worksheets.Item("About")
$wsEia   = $wb.Worksheets.Item("EIA Table 1")
$wsEoc   = $wb.Worksheets.Item("EoCEDwEC")

# ============================================================================
# 1. "About" sheet - new source citation, 3-year wording, new fuel notes
# ============================================================================

$wsAbout.Range("B4").Value = 2014
$wsAbout.Range("B5").Value = "Price Elasticities for Energy Use in Buildings of the United States"
$wsAbout.Range("B6").Value = "https://www.eia.gov/analysis/studies/buildings/energyuse/pdf/price_elasticities.pdf"
$wsAbout.Range("B7").Value = "Appendix"

# repoint the hyperlink on B6 at the new source URL
$wsAbout.Hyperlinks.Delete() | Out-Null
$wsAbout.Hyperlinks.Add($wsAbout.Range("B6"), "https://www.eia.gov/analysis/studies/buildings/energyuse/pdf/price_elasticities.pdf") | Out-Null

$wsAbout.Range("A10").Value = "We use same-price, long-run elasticities minus the 3-year short-run elasticities."
$wsAbout.Range("A11").Value = "We calculate it this way because we assume that 3-year elasticities primarily reflect behavior"
$wsAbout.Range("A14").Value = "all timescales.  So, the portion of the long-run elasticitiy represented by the 3-year elasticity"

$wsAbout.Range("A28").Value = "We assume kerosene and fuel oil-burning equipment is similar to diesel-burning equipment."
$wsAbout.Range("A29").Value = "We assume LPG/propane/butane-burning equipment is similar to natural gas-burning equipment."
$wsAbout.Range("A30").Value = "We assume hydrogen-using equipment is similar to electricity-using equipment (as it may"
$wsAbout.Range("A31").Value = "contain fuel cells that produce electricity from hydrogen)."

# ============================================================================
# 2. "EIA Table 1" sheet - refreshed Residential / Commercial elasticities
# ============================================================================

# Residential block (row 7-9)
$wsEia.Range("B7").Value = -0.12
$wsEia.Range("C7").Value = -0.21
$wsEia.Range("D7").Value = -0.25
$wsEia.Range("E7").Value = -0.28000000000000003
$wsEia.Range("F7").Value = 0
$wsEia.Range("G7").Value = 0

$wsEia.Range("B8").Value = -0.07
$wsEia.Range("C8").Value = -0.13
$wsEia.Range("D8").Value = -0.15
$wsEia.Range("E8").Value = 0.03
$wsEia.Range("F8").Value = -0.21
$wsEia.Range("G8").Value = 0

$wsEia.Range("B9").Value = -0.07
$wsEia.Range("C9").Value = -0.12
$wsEia.Range("D9").Value = -0.14000000000000001
$wsEia.Range("E9").Value = 0
$wsEia.Range("F9").Value = 0
$wsEia.Range("G9").Value = -0.22

# Commercial block (row 14-16)
$wsEia.Range("B14").Value = -0.11
$wsEia.Range("C14").Value = -0.18
$wsEia.Range("D14").Value = -0.22
$wsEia.Range("E14").Value = -0.33
$wsEia.Range("F14").Value = 0.09
$wsEia.Range("G14").Value = 0

$wsEia.Range("B15").Value = -0.15
$wsEia.Range("C15").Value = -0.25
$wsEia.Range("D15").Value = -0.3
$wsEia.Range("E15").Value = 0.15
$wsEia.Range("F15").Value = -0.57999999999999996
$wsEia.Range("G15").Value = 0.02

$wsEia.Range("B16").Value = -0.14000000000000001
$wsEia.Range("C16").Value = -0.24
$wsEia.Range("D16").Value = -0.28999999999999998
$wsEia.Range("E16").Value = 0
$wsEia.Range("F16").Value = 0.05
$wsEia.Range("G16").Value = -0.42

# ============================================================================
# 3. "EoCEDwEC" sheet - header relabel, 3-year baseline, new fuel rows
# ============================================================================

# Header: A1 becomes a wrapped, bold title; widen column A to fit it
$wsEoc.Range("A1").Value = "Elasticity by Fuel (dimensionless)"
$wsEoc.Range("A1").Font.Bold = $true
$wsEoc.Range("A1").WrapText = $true
$wsEoc.Rows.Item(1).RowHeight = 30
$wsEoc.Columns.Item(1).ColumnWidth = 23.3

# Switch the long-run-minus-short-run baseline from column B (1-year) to
# column D (3-year) of "EIA Table 1"
$wsEoc.Range("B2").Formula = "='EIA Table 1'!E7-'EIA Table 1'!D7"
$wsEoc.Range("D2").Formula = "='EIA Table 1'!E14-'EIA Table 1'!D14"
$wsEoc.Range("B4").Formula = "='EIA Table 1'!F8-'EIA Table 1'!D8"
$wsEoc.Range("D4").Formula = "='EIA Table 1'!F15-'EIA Table 1'!D15"
$wsEoc.Range("B5").Formula = "='EIA Table 1'!G9-'EIA Table 1'!D9"
$wsEoc.Range("D5").Formula = "='EIA Table 1'!G16-'EIA Table 1'!D16"

# New fuel rows: kerosene / heavy fuel oil behave like petroleum diesel (row 5);
# LPG/propane/butane behaves like natural gas (row 4); hydrogen behaves like
# electricity (row 2). All four new rows get the same orange fill used for
# the other "assumed similar to" rows (coal/petroleum diesel).
$wsEoc.Range("A8").Value = "kerosene"
$wsEoc.Range("B8").Formula = "=B5"
$wsEoc.Range("C8").Formula = "=C5"
$wsEoc.Range("D8").Formula = "=D5"

$wsEoc.Range("A9").Value = "heavy or residual fuel oil"
$wsEoc.Range("B9").Formula = "=B5"
$wsEoc.Range("C9").Formula = "=C5"
$wsEoc.Range("D9").Formula = "=D5"

$wsEoc.Range("A10").Value = "LPG propane or butane"
$wsEoc.Range("B10").Formula = "=B4"
$wsEoc.Range("C10").Formula = "=C4"
$wsEoc.Range("D10").Formula = "=D4"

$wsEoc.Range("A11").Value = "hydrogen"
$wsEoc.Range("B11").Formula = "=B2"
$wsEoc.Range("C11").Formula = "=C2"
$wsEoc.Range("D11").Formula = "=D2"

$wsEoc.Range("B8:D11").Interior.Color = 49407

# ============================================================================
# 4. Selections / active sheet, matching the saved view in the edited file
# ============================================================================

$wsAbout.Range("A27").Select() | Out-Null
$wsEia.Range("E17").Select() | Out-Null
$wsEoc.Range("H29").Select() | Out-Null
$wsEoc.Activate() | Out-Null
